$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update unit label from "mg/l" to "mg/mL" (case-insensitive unit parsing fix)
$ws.Range("B2").Value = "mg/mL"

# Update selection / active cell on sheet1
$ws.Range("D9").Select()

# Widen column B to fit the new, longer unit label ("mg/mL")
$ws.Columns.Item(2).ColumnWidth = 13.43
